# Added ifo GDP component analysis preprocessing:
# populate additional diagonal values in the matched-errors matrix
# (one extra column of data per row, continuing the staircase pattern)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: update existing value
$ws.Range("K14").Value = 0.5403464745801891

# Row 15: update existing value, add new one
$ws.Range("J15").Value = 0.4969491838668565
$ws.Range("K15").Value = 0.2970525035592049

# Row 16: update existing value, add new one
$ws.Range("I16").Value = 0.5011245128056051
$ws.Range("J16").Value = 0.2858677898194339

# Row 17: update existing value, add new one
$ws.Range("H17").Value = 0.4852787037784192
$ws.Range("I17").Value = 0.2775335613519331

# Row 18: update existing value, add new one
$ws.Range("G18").Value = 0.4539510573947921
$ws.Range("H18").Value = 0.2743085116504074

# Row 19: update existing value, add new one
$ws.Range("F19").Value = 0.4663391832225094
$ws.Range("G19").Value = 0.2534447081011285

# Row 20: update existing value, add new one
$ws.Range("E20").Value = 0.4814444548743619
$ws.Range("F20").Value = 0.2766837437271186

# Row 21: update existing value, add new one
$ws.Range("D21").Value = 0.4184715358843989
$ws.Range("E21").Value = 0.2867219094086165

# Row 22: update existing value, add new one
$ws.Range("C22").Value = 0.5177895860664353
$ws.Range("D22").Value = 0.1751453671933744

# Row 23: update existing value, add new one
$ws.Range("B23").Value = 0.5618492773058843
$ws.Range("C23").Value = 0.1965658720679752

# Row 24: add new value
$ws.Range("B24").Value = 0.4328090033804217
